$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 699465.0600000001
$ws.Range("I98").Value = 699465.0600000001
$ws.Range("K98").Value = 699465.0600000001
$ws.Range("M98").Value = -697967.0600000001
$ws.Range("H122").Value = 699465.0600000001
$ws.Range("I122").Value = 699465.0600000001
$ws.Range("K122").Value = 2098395.18
$ws.Range("M122").Value = -2095945.18
$ws.Range("H141").Value = 2754.55
$ws.Range("I141").Value = 1971.7273
$ws.Range("J141").Value = 6445
$ws.Range("K141").Value = 5915.1819
$ws.Range("L141").Value = 19335
$ws.Range("M141").Value = -735.1818999999996
$ws.Range("N141").Value = -29695

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31177.59
$ws.Range("I32").Value = 4926.8237
$ws.Range("J32").Value = 209682.8
$ws.Range("K32").Value = 4926.8237
$ws.Range("L32").Value = 209682.8
$ws.Range("M32").Value = -4639.8237
$ws.Range("N32").Value = -210256.8
$ws.Range("H45").Value = 1356
$ws.Range("I45").Value = 1212
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 1212
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -835
$ws.Range("N45").Value = -2254
$ws.Range("H74").Value = 6528.5835
$ws.Range("I74").Value = 987.06665
$ws.Range("J74").Value = 15764.444
$ws.Range("K74").Value = 987.06665
$ws.Range("L74").Value = 15764.444
$ws.Range("M74").Value = -113.06665
$ws.Range("N74").Value = -17512.444
$ws.Range("H77").Value = 6528.5835
$ws.Range("I77").Value = 987.06665
$ws.Range("J77").Value = 15764.444
$ws.Range("K77").Value = 4935.33325
$ws.Range("L77").Value = 78822.22
$ws.Range("M77").Value = -567.3332499999997
$ws.Range("N77").Value = -87558.22
$ws.Range("H122").Value = 2633
$ws.Range("I122").Value = 1232.875
$ws.Range("K122").Value = 3698.625
$ws.Range("M122").Value = -1248.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 956.75757
$ws.Range("I20").Value = 845
$ws.Range("J20").Value = 1213.8
$ws.Range("K20").Value = 845
$ws.Range("L20").Value = 1213.8
$ws.Range("M20").Value = -598
$ws.Range("N20").Value = -1707.8
$ws.Range("H86").Value = 15601.857
$ws.Range("I86").Value = 1568.6666
$ws.Range("J86").Value = 26126.75
$ws.Range("K86").Value = 1568.6666
$ws.Range("L86").Value = 26126.75
$ws.Range("M86").Value = -445.6666
$ws.Range("N86").Value = -28372.75
$ws.Range("H89").Value = 15601.857
$ws.Range("I89").Value = 1568.6666
$ws.Range("J89").Value = 26126.75
$ws.Range("K89").Value = 7843.333000000001
$ws.Range("L89").Value = 130633.75
$ws.Range("M89").Value = -2227.333000000001
$ws.Range("N89").Value = -141865.75
$ws.Range("H94").Value = 1337.5
$ws.Range("I94").Value = 1337.5
$ws.Range("K94").Value = 1337.5
$ws.Range("M94").Value = -886.5
$ws.Range("H105").Value = 220081.9
$ws.Range("I105").Value = 2489.1177
$ws.Range("J105").Value = 836594.8
$ws.Range("K105").Value = 2489.1177
$ws.Range("L105").Value = 836594.8
$ws.Range("M105").Value = -742.1176999999998
$ws.Range("N105").Value = -840088.8
$ws.Range("H107").Value = 3055.5
$ws.Range("I107").Value = 3055.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3055.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1135.5
$ws.Range("N107").ClearContents()
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H111").Value = 25000
$ws.Range("J111").Value = 25000
$ws.Range("L111").Value = 25000
$ws.Range("N111").Value = -33180
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1071.4286
$ws.Range("I16").Value = 1120
$ws.Range("J16").Value = 950
$ws.Range("K16").Value = 1120
$ws.Range("L16").Value = 950
$ws.Range("M16").Value = -833
$ws.Range("N16").Value = -1524
$ws.Range("H31").Value = 3091.6033
$ws.Range("I31").Value = 1185.1428
$ws.Range("K31").Value = 1185.1428
$ws.Range("M31").Value = -890.1428000000001
$ws.Range("H34").Value = 3091.6033
$ws.Range("I34").Value = 1185.1428
$ws.Range("K34").Value = 1185.1428
$ws.Range("M34").Value = -983.1428000000001
$ws.Range("H113").Value = 1071.4286
$ws.Range("I113").Value = 1120
$ws.Range("J113").Value = 950
$ws.Range("K113").Value = 1120
$ws.Range("L113").Value = 950
$ws.Range("M113").Value = 1050
$ws.Range("N113").Value = -5290

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 461.64
$ws.Range("J107").Value = 397.7143
$ws.Range("L107").Value = 1193.1429
$ws.Range("N107").Value = -5033.1429
$ws.Range("H122").Value = 1173.24
$ws.Range("I122").Value = 196
$ws.Range("K122").Value = 1764
$ws.Range("M122").Value = 686

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3372.0557
$ws.Range("I102").Value = 3413.5334
$ws.Range("J102").Value = 3164.6667
$ws.Range("K102").Value = 3413.5334
$ws.Range("L102").Value = 3164.6667
$ws.Range("M102").Value = -1791.5334
$ws.Range("N102").Value = -6408.6667
$ws.Range("H122").Value = 696088.4
$ws.Range("I122").Value = 926972.3
$ws.Range("J122").Value = 3436.5
$ws.Range("K122").Value = 2780916.9
$ws.Range("L122").Value = 10309.5
$ws.Range("M122").Value = -2778466.9
$ws.Range("N122").Value = -15209.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2797
$ws.Range("I7").Value = 2200.4
$ws.Range("J7").Value = 3068.182
$ws.Range("K7").Value = 2200.4
$ws.Range("L7").Value = 3068.182
$ws.Range("M7").Value = -2088.4
$ws.Range("N7").Value = -3292.182
$ws.Range("H40").Value = 2439
$ws.Range("I40").Value = 1333.1111
$ws.Range("J40").Value = 2991.9443
$ws.Range("K40").Value = 1333.1111
$ws.Range("L40").Value = 2991.9443
$ws.Range("M40").Value = -1197.1111
$ws.Range("N40").Value = -3263.9443
$ws.Range("H126").Value = 2797
$ws.Range("I126").Value = 2200.4
$ws.Range("J126").Value = 3068.182
$ws.Range("K126").Value = 6601.200000000001
$ws.Range("L126").Value = 9204.545999999998
$ws.Range("M126").Value = -4131.200000000001
$ws.Range("N126").Value = -14144.546
$ws.Range("H133").Value = 47062.92
$ws.Range("J133").Value = 47062.92
$ws.Range("L133").Value = 47062.92
$ws.Range("N133").Value = -52122.92
